$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base input value (B1) which drives the dependent formulas
# in row 8 (B8:M8) and, transitively, row 10 (B10:M10).
$ws.Range("B1").Value = 1

# Update the active cell selection to match the saved view state.
$ws.Range("D15").Select()
